$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Revise "Data" values for existing years (rows 2-142) ---
# Values are written as a TEXT() formula first (so Excel stores them as the
# literal text strings the source file uses, not as numbers), then the
# formulas are converted to plain values via Copy / PasteSpecial -Values so no
# formula or extra cell style is left behind.
$ws.Range("E2").Formula = "=TEXT(1423,""0.############"")"
$ws.Range("E32").Formula = "=TEXT(2206,""0.############"")"
$ws.Range("E33").Formula = "=TEXT(2429,""0.############"")"
$ws.Range("E34").Formula = "=TEXT(2297,""0.############"")"
$ws.Range("E35").Formula = "=TEXT(1975,""0.############"")"
$ws.Range("E36").Formula = "=TEXT(1910,""0.############"")"
$ws.Range("E37").Formula = "=TEXT(2114,""0.############"")"
$ws.Range("E38").Formula = "=TEXT(2375,""0.############"")"
$ws.Range("E39").Formula = "=TEXT(2166,""0.############"")"
$ws.Range("E40").Formula = "=TEXT(2488,""0.############"")"
$ws.Range("E41").Formula = "=TEXT(2515,""0.############"")"
$ws.Range("E42").Formula = "=TEXT(2732,""0.############"")"
$ws.Range("E43").Formula = "=TEXT(2447,""0.############"")"
$ws.Range("E44").Formula = "=TEXT(2094,""0.############"")"
$ws.Range("E45").Formula = "=TEXT(1812,""0.############"")"
$ws.Range("E46").Formula = "=TEXT(2166,""0.############"")"
$ws.Range("E47").Formula = "=TEXT(2209,""0.############"")"
$ws.Range("E48").Formula = "=TEXT(1827,""0.############"")"
$ws.Range("E49").Formula = "=TEXT(2530,""0.############"")"
$ws.Range("E50").Formula = "=TEXT(1978,""0.############"")"
$ws.Range("E51").Formula = "=TEXT(2244,""0.############"")"
$ws.Range("E52").Formula = "=TEXT(1932,""0.############"")"
$ws.Range("E53").Formula = "=TEXT(2182,""0.############"")"
$ws.Range("E54").Formula = "=TEXT(2884,""0.############"")"
$ws.Range("E55").Formula = "=TEXT(2952,""0.############"")"
$ws.Range("E56").Formula = "=TEXT(2888,""0.############"")"
$ws.Range("E57").Formula = "=TEXT(3255,""0.############"")"
$ws.Range("E58").Formula = "=TEXT(2686,""0.############"")"
$ws.Range("E59").Formula = "=TEXT(1986,""0.############"")"
$ws.Range("E60").Formula = "=TEXT(2361,""0.############"")"
$ws.Range("E61").Formula = "=TEXT(3323,""0.############"")"
$ws.Range("E62").Formula = "=TEXT(3049,""0.############"")"
$ws.Range("E63").Formula = "=TEXT(2246,""0.############"")"
$ws.Range("E64").Formula = "=TEXT(1954,""0.############"")"
$ws.Range("E65").Formula = "=TEXT(2766,""0.############"")"
$ws.Range("E66").Formula = "=TEXT(3263,""0.############"")"
$ws.Range("E67").Formula = "=TEXT(3676,""0.############"")"
$ws.Range("E68").Formula = "=TEXT(3652,""0.############"")"
$ws.Range("E69").Formula = "=TEXT(3967,""0.############"")"
$ws.Range("E70").Formula = "=TEXT(3363,""0.############"")"
$ws.Range("E71").Formula = "=TEXT(3912,""0.############"")"
$ws.Range("E82").Formula = "=TEXT(3572,""0.############"")"
$ws.Range("E83").Formula = "=TEXT(3368,""0.############"")"
$ws.Range("E84").Formula = "=TEXT(2342,""0.############"")"
$ws.Range("E85").Formula = "=TEXT(2794,""0.############"")"
$ws.Range("E86").Formula = "=TEXT(3202,""0.############"")"
$ws.Range("E87").Formula = "=TEXT(2946,""0.############"")"
$ws.Range("E88").Formula = "=TEXT(3241,""0.############"")"
$ws.Range("E89").Formula = "=TEXT(2874,""0.############"")"
$ws.Range("E90").Formula = "=TEXT(2880,""0.############"")"
$ws.Range("E91").Formula = "=TEXT(3239,""0.############"")"
$ws.Range("E92").Formula = "=TEXT(3464,""0.############"")"
$ws.Range("E93").Formula = "=TEXT(3623,""0.############"")"
$ws.Range("E94").Formula = "=TEXT(3775,""0.############"")"
$ws.Range("E95").Formula = "=TEXT(4049,""0.############"")"
$ws.Range("E96").Formula = "=TEXT(3800,""0.############"")"
$ws.Range("E97").Formula = "=TEXT(3991,""0.############"")"
$ws.Range("E98").Formula = "=TEXT(4317,""0.############"")"
$ws.Range("E99").Formula = "=TEXT(4741,""0.############"")"
$ws.Range("E100").Formula = "=TEXT(5294,""0.############"")"
$ws.Range("E101").Formula = "=TEXT(5931,""0.############"")"
$ws.Range("E102").Formula = "=TEXT(6650,""0.############"")"
$ws.Range("E103").Formula = "=TEXT(7318,""0.############"")"
$ws.Range("E104").Formula = "=TEXT(8155,""0.############"")"
$ws.Range("E105").Formula = "=TEXT(8894,""0.############"")"
$ws.Range("E106").Formula = "=TEXT(9314,""0.############"")"
$ws.Range("E107").Formula = "=TEXT(9602,""0.############"")"
$ws.Range("E108").Formula = "=TEXT(10178,""0.############"")"
$ws.Range("E109").Formula = "=TEXT(10790,""0.############"")"
$ws.Range("E110").Formula = "=TEXT(11588,""0.############"")"
$ws.Range("E111").Formula = "=TEXT(12521,""0.############"")"
$ws.Range("E112").Formula = "=TEXT(13601,""0.############"")"
$ws.Range("E113").Formula = "=TEXT(14335,""0.############"")"
$ws.Range("E114").Formula = "=TEXT(14687,""0.############"")"
$ws.Range("E115").Formula = "=TEXT(15720,""0.############"")"
$ws.Range("E116").Formula = "=TEXT(16767,""0.############"")"
$ws.Range("E117").Formula = "=TEXT(16611,""0.############"")"
$ws.Range("E118").Formula = "=TEXT(16831,""0.############"")"
$ws.Range("E119").Formula = "=TEXT(18345,""0.############"")"
$ws.Range("E120").Formula = "=TEXT(19853,""0.############"")"
$ws.Range("E121").Formula = "=TEXT(21221,""0.############"")"
$ws.Range("E122").Formula = "=TEXT(22666,""0.############"")"
$ws.Range("E123").Formula = "=TEXT(23636.6826875653,""0.############"")"
$ws.Range("E124").Formula = "=TEXT(24775.5711919494,""0.############"")"
$ws.Range("E125").Formula = "=TEXT(27073.714552744,""0.############"")"
$ws.Range("E126").Formula = "=TEXT(29454.6124611513,""0.############"")"
$ws.Range("E127").Formula = "=TEXT(30957.2871651082,""0.############"")"
$ws.Range("E128").Formula = "=TEXT(32724.4875328756,""0.############"")"
$ws.Range("E129").Formula = "=TEXT(34868.0659692042,""0.############"")"
$ws.Range("E130").Formula = "=TEXT(33590.4819258856,""0.############"")"
$ws.Range("E131").Formula = "=TEXT(35156.1182739866,""0.############"")"
$ws.Range("E132").Formula = "=TEXT(37772.7596926971,""0.############"")"
$ws.Range("E133").Formula = "=TEXT(36864.5751243973,""0.############"")"
$ws.Range("E134").Formula = "=TEXT(37829.9043369914,""0.############"")"
$ws.Range("E135").Formula = "=TEXT(38950.7459376796,""0.############"")"
$ws.Range("E136").Formula = "=TEXT(42111.1662234909,""0.############"")"
$ws.Range("E137").Formula = "=TEXT(44698.8385991564,""0.############"")"
$ws.Range("E138").Formula = "=TEXT(48073.5880877817,""0.############"")"
$ws.Range("E139").Formula = "=TEXT(51844.0617998722,""0.############"")"
$ws.Range("E140").Formula = "=TEXT(52180.1731759185,""0.############"")"
$ws.Range("E141").Formula = "=TEXT(51344.4780400612,""0.############"")"
$ws.Range("E142").Formula = "=TEXT(58612.7310490197,""0.############"")"

# --- Append new rows for years 2011-2016 ---
$ws.Range("A143").Value = 702
$ws.Range("B143").Formula = "=""Singapore"""
$ws.Range("C143").Formula = "=""GDP per Capita"""
$ws.Range("D143").Value = 2011
$ws.Range("E143").Formula = "=TEXT(61672,""0.############"")"
$ws.Range("A144").Value = 702
$ws.Range("B144").Formula = "=""Singapore"""
$ws.Range("C144").Formula = "=""GDP per Capita"""
$ws.Range("D144").Value = 2012
$ws.Range("E144").Formula = "=TEXT(62783,""0.############"")"
$ws.Range("A145").Value = 702
$ws.Range("B145").Formula = "=""Singapore"""
$ws.Range("C145").Formula = "=""GDP per Capita"""
$ws.Range("D145").Value = 2013
$ws.Range("E145").Formula = "=TEXT(64633,""0.############"")"
$ws.Range("A146").Value = 702
$ws.Range("B146").Formula = "=""Singapore"""
$ws.Range("C146").Formula = "=""GDP per Capita"""
$ws.Range("D146").Value = 2014
$ws.Range("E146").Formula = "=TEXT(65655,""0.############"")"
$ws.Range("A147").Value = 702
$ws.Range("B147").Formula = "=""Singapore"""
$ws.Range("C147").Formula = "=""GDP per Capita"""
$ws.Range("D147").Value = 2015
$ws.Range("E147").Formula = "=TEXT(65660,""0.############"")"
$ws.Range("A148").Value = 702
$ws.Range("B148").Formula = "=""Singapore"""
$ws.Range("C148").Formula = "=""GDP per Capita"""
$ws.Range("D148").Value = 2016
$ws.Range("E148").Formula = "=TEXT(65729,""0.############"")"

# --- Convert all the Formula-driven cells above to literal values ---
$textRng = $ws.Range("E2:E148")
$textRng.Copy()
$textRng.PasteSpecial(-4163)
$strRng = $ws.Range("B143:C148")
$strRng.Copy()
$strRng.PasteSpecial(-4163)
$excel.CutCopyMode = 0
